$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.111.12"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "2.544.93"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.00"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.66"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("D9").Value = "2.543.79"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("E10").Value = "  +0.38%  "
$ws.Range("E11").Value = "  +1.78%  "
$ws.Range("E12").Value = "  -1.85%  "
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.54"
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "3.002.29"
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("D17").Value = "67.855.73"
$ws.Range("E17").Value = "  +1.16%  "
$ws.Range("E18").Value = "  +137.34%  "
$ws.Range("D19").Value = "2.542.53"
$ws.Range("E19").Value = "  +1.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.81"
$ws.Range("E20").Value = "  +3.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.01"
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "370.09"
$ws.Range("E22").Value = "  +4.15%  "
$ws.Range("E23").Value = "  -0.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.60"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.75"
$ws.Range("E26").Value = "  -2.88%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.99"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("E29").Value = "  -1.58%  "
$ws.Range("D30").Value = "0.0₃0976"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.48"
$ws.Range("E31").Value = "  +3.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "543.65"
$ws.Range("E32").Value = "  +0.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.33"
$ws.Range("E33").Value = "  -0.64%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.88"
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.53"
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("E38").Value = "  -1.53%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.19"
$ws.Range("E39").Value = "  +2.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.62"
$ws.Range("E40").Value = "  +0.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.17"
$ws.Range("E41").Value = "  +0.61%  "
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.58"
$ws.Range("E44").Value = "  +2.89%  "
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.24"
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("E47").Value = "  +4.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "147.93"
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.73"
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("E51").Value = "  +1.38%  "
